# Speaker bio fixes on slide 2 ("Speakers") + repositioning of the second
# "LesFurets.wmf" picture, per the "fixed speaker bio + slides PDF" commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# Shape 3 ("Espace réservé du contenu 4"): fix the first speaker's bio.
#   "Developpeur" + " de plus de 30 ans"  ->  "Développeur " + "de plus de 30 ans"
# ---------------------------------------------------------------------
$sh3 = $s.Shapes.Item(3)

$t = $sh3.TextFrame.TextRange.Text
$idx = $t.IndexOf("Developpeur")
$sub = $sh3.TextFrame.TextRange.Characters($idx + 1, "Developpeur".Length)
$sub.Text = "Développeur "

$t = $sh3.TextFrame.TextRange.Text
$idx = $t.IndexOf(" de plus de 30 ans")
$sub = $sh3.TextFrame.TextRange.Characters($idx + 1, " de plus de 30 ans".Length)
$sub.Text = "de plus de 30 ans"

# ---------------------------------------------------------------------
# Shape 5 ("Espace réservé du contenu 6"): fix the second speaker's bio.
#   "Dévelopeur" + " java depuis 1999"  ->  "Développeur " + "java depuis 1999"
#   "Architecte pour" -> "Architecte " + "pour", followed by two new blank lines
# ---------------------------------------------------------------------
$sh5 = $s.Shapes.Item(5)

$t = $sh5.TextFrame.TextRange.Text
$idx = $t.IndexOf("Dévelopeur")
$sub = $sh5.TextFrame.TextRange.Characters($idx + 1, "Dévelopeur".Length)
$sub.Text = "Développeur "

$t = $sh5.TextFrame.TextRange.Text
$idx = $t.IndexOf(" java depuis 1999")
$sub = $sh5.TextFrame.TextRange.Characters($idx + 1, " java depuis 1999".Length)
$sub.Text = "java depuis 1999"

# Insert two new blank paragraphs right after the "Architecte pour" paragraph,
# by inserting them before the following paragraph ("ILOG - IBM", #3).
$tr5 = $sh5.TextFrame.TextRange
$paraIlog = $tr5.Paragraphs(3, 1)
$paraIlog.InsertBefore("`r`r")

# Now split "Architecte pour" into two runs: "Architecte " and "pour"
$t = $sh5.TextFrame.TextRange.Text
$idx = $t.IndexOf("Architecte pour")
$sub = $sh5.TextFrame.TextRange.Characters($idx + 1, "Architecte pour".Length)
$sub.Text = "Architecte "

$t = $sh5.TextFrame.TextRange.Text
$idx = $t.IndexOf("Architecte ")
$sub = $sh5.TextFrame.TextRange.Characters($idx + 1, "Architecte ".Length)
$newRun = $sub.InsertAfter("pour")
$newRun.Font.Bold = $true

# ---------------------------------------------------------------------
# Reposition the second "LesFurets.wmf" picture (Image 12).
# ---------------------------------------------------------------------
$pic = $s.Shapes.Item(8)
$pic.Left = 5562600 / 12700
$pic.Top = 3124200 / 12700
